$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Chilean Primera Division) is removed entirely; Excel shifts rows 3 & 4 up to 2 & 3.
$ws.Rows.Item(2).Delete()

# Prevent the "Date" column (B) text from being auto-converted to a date serial number
# by temporarily marking it as Text, writing the literal string, then clearing the format
# back off so the cell ends up with no explicit style (matching the rest of the data rows).
$ws.Range("B2:B3").NumberFormat = "@"

# Row 2 now holds what used to be row 3: Argentinian Primera Division / Racing Club vs Tigre,
# but with the updated odds below.
$ws.Range("A2").Value2 = 'Argentinian Primera Division'
$ws.Range("B2").Value2 = '2025-12-01'
$ws.Range("C2").Value2 = '21:30:00'
$ws.Range("D2").Value2 = 'Racing Club'
$ws.Range("E2").Value2 = 'Tigre'
$ws.Range("F2").Value2 = 3.45
$ws.Range("G2").Value2 = 3.55
$ws.Range("H2").Value2 = 12
$ws.Range("I2").Value2 = 13
$ws.Range("J2").Value2 = 1.56
$ws.Range("K2").Value2 = 1.57
$ws.Range("L2").Value2 = 0
$ws.Range("M2").Value2 = 2.38
$ws.Range("N2").Value2 = 1.1
$ws.Range("O2").Value2 = 10
$ws.Range("P2").Value2 = 1.01
$ws.Range("Q2").Value2 = 65
$ws.Range("R2").Value2 = 1.01
$ws.Range("S2").Value2 = 420
$ws.Range("T2").Value2 = 21
$ws.Range("U2").Value2 = 1.04
$ws.Range("V2").Value2 = 1.09
$ws.Range("W2").Value2 = 1.4
$ws.Range("X2").Value2 = 1.73
$ws.Range("Y2").Value2 = 13.5
$ws.Range("Z2").Value2 = 1000
$ws.Range("AA2").Value2 = 1000
$ws.Range("AB2").Value2 = 4.2
$ws.Range("AC2").Value2 = 30
$ws.Range("AD2").Value2 = 1000
$ws.Range("AE2").Value2 = 1000
$ws.Range("AF2").Value2 = 28
$ws.Range("AG2").Value2 = 1000
$ws.Range("AH2").Value2 = 1000
$ws.Range("AI2").Value2 = 1000
$ws.Range("AJ2").Value2 = 370
$ws.Range("AK2").Value2 = 1000
$ws.Range("AL2").Value2 = 1000
$ws.Range("AM2").Value2 = 1000
$ws.Range("AN2").Value2 = 1000
$ws.Range("AO2").Value2 = 1000

# Row 3 now holds what used to be row 4: Colombian Primera A / Ind Medellin vs America de Cali S.A,
# but with the updated odds below.
$ws.Range("A3").Value2 = 'Colombian Primera A'
$ws.Range("B3").Value2 = '2025-12-01'
$ws.Range("C3").Value2 = '22:00:00'
$ws.Range("D3").Value2 = 'Ind Medellin'
$ws.Range("E3").Value2 = 'America de Cali S.A'
$ws.Range("F3").Value2 = 2.24
$ws.Range("G3").Value2 = 2.28
$ws.Range("H3").Value2 = 5.5
$ws.Range("I3").Value2 = 5.8
$ws.Range("J3").Value2 = 2.58
$ws.Range("K3").Value2 = 2.66
$ws.Range("L3").Value2 = 0
$ws.Range("M3").Value2 = 0
$ws.Range("N3").Value2 = 0
$ws.Range("O3").Value2 = 0
$ws.Range("P3").Value2 = 4
$ws.Range("Q3").Value2 = 1.32
$ws.Range("R3").Value2 = 1.64
$ws.Range("S3").Value2 = 2.48
$ws.Range("T3").Value2 = 0
$ws.Range("U3").Value2 = 0
$ws.Range("V3").Value2 = 1.21
$ws.Range("W3").Value2 = 1.77
$ws.Range("X3").Value2 = 1000
$ws.Range("Y3").Value2 = 1000
$ws.Range("Z3").Value2 = 1000
$ws.Range("AA3").Value2 = 1000
$ws.Range("AB3").Value2 = 1000
$ws.Range("AC3").Value2 = 4.2
$ws.Range("AD3").Value2 = 9.2
$ws.Range("AE3").Value2 = 50
$ws.Range("AF3").Value2 = 1000
$ws.Range("AG3").Value2 = 4.5
$ws.Range("AH3").Value2 = 8.2
$ws.Range("AI3").Value2 = 55
$ws.Range("AJ3").Value2 = 1000
$ws.Range("AK3").Value2 = 11.5
$ws.Range("AL3").Value2 = 21
$ws.Range("AM3").Value2 = 100
$ws.Range("AN3").Value2 = 27
$ws.Range("AO3").Value2 = 150

# Drop the temporary text formatting on column B so the cells carry no explicit style,
# matching the rest of the (unstyled) data rows.
$ws.Range("B2:B3").ClearFormats()
